$wb = $excel.ActiveWorkbook

# zh-cn sheet: update the "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# values for the bd5c11a7... row (row 4).
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D4").Value = "2016-01-25 06:22:55"
$wsZhCn.Range("G4").Value = "2016-01-25 06:23:39"

# de-de sheet: same row/columns.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D4").Value = "2016-01-25 06:23:06"
$wsDeDe.Range("G4").Value = "2016-01-25 06:23:58"
